$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.172.31"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.59"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.72"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("E9").Value = "  +3.48%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.44"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "3.030.81"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "63.142.37"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "2.583.23"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.59"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.42"
$ws.Range("E20").Value = "  +3.08%  "
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -3.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.89"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").Value = "2.697.23"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +12.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.48"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("E32").Value = "  +8.07%  "
$ws.Range("D33").Value = "0.0₃0828"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "466.10"
$ws.Range("E34").Value = "  +14.52%  "
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "175.76"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.22"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.56"
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.86"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.03"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.613"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0979"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.75"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("E51").Value = "  -0.08%  "
